$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data per diff
# Row 2
$ws.Range("D2").Value = "65.914.23"
$ws.Range("E2").Value = "  +2.02%  "
# Row 3
$ws.Range("D3").Value = "3.381.20"
$ws.Range("E3").Value = "  +1.49%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "183.00"
$ws.Range("E5").Value = "  +0.88%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "541.45"
$ws.Range("E6").Value = "  +1.74%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").Value = "  +0.05%  "
# Row 8
$ws.Range("D8").Value = "3.377.14"
$ws.Range("E8").Value = "  +1.55%  "
# Row 9
$ws.Range("E9").Value = "  -0.08%  "
# Row 10
$ws.Range("E10").Value = "  +2.76%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.25"
$ws.Range("E11").Value = "  -5.80%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.145"
$ws.Range("E12").Value = "  +7.93%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("E13").Value = "  +2.88%  "
# Row 14
$ws.Range("E14").Value = "  +1.22%  "
# Row 15
$ws.Range("D15").Value = "3.911.62"
$ws.Range("E15").Value = "  +1.13%  "
# Row 16
$ws.Range("E16").Value = "  +2.03%  "
# Row 17
$ws.Range("D17").Value = "3.373.46"
$ws.Range("E17").Value = "  +1.14%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.08"
$ws.Range("E18").Value = "  +2.74%  "
# Row 19
$ws.Range("D19").Value = "66.013.37"
$ws.Range("E19").Value = "  +2.16%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.55"
$ws.Range("E20").Value = "  +2.93%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.987"
$ws.Range("E21").Value = "  +2.05%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "402.53"
$ws.Range("E22").Value = "  +6.78%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.88"
$ws.Range("E23").Value = "  +5.44%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.22"
$ws.Range("E24").Value = "  +6.57%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.86"
$ws.Range("E25").Value = "  +3.21%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.83"
$ws.Range("E26").Value = "  +0.04%  "
# Row 27
$ws.Range("E27").Value = "  +5.50%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.12"
$ws.Range("E28").Value = "  +0.37%  "
# Row 29
$ws.Range("E29").Value = "  +1.03%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.51"
$ws.Range("E30").Value = "  +0.79%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.79"
$ws.Range("E31").Value = "  +2.00%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "664.32"
$ws.Range("E32").Value = "  +1.62%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.78"
$ws.Range("E33").Value = "  +1.06%  "
# Row 34
$ws.Range("E34").Value = "  +1.11%  "
# Row 35
$ws.Range("E35").Value = "  +1.73%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.10"
$ws.Range("E36").Value = "  -2.69%  "
# Row 37
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0817"
$ws.Range("E37").Value = "  +14.66%  "
# Row 38
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.15"
$ws.Range("E38").Value = "  +2.84%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.06%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.397"
$ws.Range("E40").Value = "  +0.82%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  +9.75%  "
# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.131"
$ws.Range("E42").Value = "  +4.30%  "
# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.29"
$ws.Range("E43").Value = "  +17.06%  "
# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.16%  "
# Row 45
$ws.Range("D45").Value = "3.014.24"
$ws.Range("E45").Value = "  +2.40%  "
# Row 46
$ws.Range("E46").Value = "  +2.70%  "
# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0415"
$ws.Range("E47").Value = "  +3.25%  "
# Row 48
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("E48").Value = "  +2.74%  "
# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.21"
$ws.Range("E49").Value = "  +6.07%  "
# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.79"
$ws.Range("E50").Value = "  +10.19%  "
# Row 51
$ws.Range("E51").Value = "  +1.29%  "
